$d = $word.ActiveDocument

# --- 1. Merge the "Articulate " / "the purpose..." runs into a single run ---
$articulateIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Articulate *purpose of your bot*") {
        $articulateIdx = $i
    }
}
if ($articulateIdx -ne -1) {
    $pr = $d.Paragraphs.Item($articulateIdx).Range
    # Range excluding the trailing paragraph mark
    $textRange = $d.Range($pr.Start, $pr.End - 1)
    $textRange.Delete()
    $textRange.InsertAfter("Articulate the purpose of your bot and take special care if your bot will support consequential use cases")
}

# --- 2. Add two new bullet paragraphs after "Bot lets user know..." ---
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Bot lets user know right off the bat*") {
        $targetIdx = $i
    }
}
if ($targetIdx -ne -1) {
    $target = $d.Paragraphs.Item($targetIdx)
    $target.Range.InsertParagraphAfter()
    $p1 = $d.Paragraphs.Item($targetIdx + 1)
    $p1.Range.InsertBefore("Some phrases were found on the internet")
    $p1.Range.ListFormat.ListLevelNumber = 1

    $p1.Range.InsertParagraphAfter()
    $p2 = $d.Paragraphs.Item($targetIdx + 2)
    $p2.Range.InsertBefore("Those demo phrases are not part of knowledgebase")
    $p2.Range.ListFormat.ListLevelNumber = 1
}
